$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of expense data (row 6) -- copy formatting from the row above
# first so the new date cell reuses the existing date style (YYYY-MM-DD).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("A6").Value = 43893
$ws.Range("B6").Value = "DEBORA"
$ws.Range("C6").Value = 23
$ws.Range("D6").Value = "Aniver"

# Update selection to match the recorded end-state
$ws.Range("E10").Select() | Out-Null
